$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("D3").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E3").Value = "['Normal']"

# Row 26
$ws.Range("D26").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['SoftwareFault']"

# Row 38
$ws.Range("D38").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E38").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'SoftwareFault']"

# Row 53
$ws.Range("D53").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E53").Value = "['Normal']"

# Row 56
$ws.Range("D56").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['Normal']"

# Row 71
$ws.Range("D71").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E71").Value = "['Normal']"

# Row 82
$ws.Range("D82").Value = "[1, 1, 1, 0, 0, 0, 0]"
$ws.Range("E82").Value = "['Normal', 'SurroundingEnvironment', 'HardwareFault']"

# Row 83
$ws.Range("D83").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E83").Value = "['Normal']"

# Row 84
$ws.Range("D84").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E84").Value = "['Normal']"
